$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# === Apply cell value changes (order matters: matches new-shared-string insertion order) ===
$v9c7 = @'
What name do I need to use for renaming columns if the table is a map.
'@
$ws.Cells.Item(9,7).Value = $v9c7

$v8c1 = @'
Set-Intersection Operation
'@
$ws.Cells.Item(8,1).Value = $v8c1

$v8c2 = @'
Tuples that exist in both relations.
'@
$ws.Cells.Item(8,2).Value = $v8c2

$v8c4 = @'
select * from Professor
INTERSECT
select * from Professor2
'@
$ws.Cells.Item(8,4).Value = $v8c4

$v8c5 = @'
professors.stream().filter(professors2::contains).collect(Collectors.toSet());
'@
$ws.Cells.Item(8,5).Value = $v8c5

$v8c6 = @'
professors.stream().filter(professors2::contains).collect(Collectors.toSet());
'@
$ws.Cells.Item(8,6).Value = $v8c6

$v8c7 = @'
Forcing using maps for this operation?
'@
$ws.Cells.Item(8,7).Value = $v8c7

$v9c1 = @'
Cartesian-product
'@
$ws.Cells.Item(9,1).Value = $v9c1

$v9c2 = @'
Join relations with different attributes.
'@
$ws.Cells.Item(9,2).Value = $v9c2

$v9c4 = @'
Professor x Department
'@
$ws.Cells.Item(9,4).Value = $v9c4

$v9c5 = ""
$ws.Cells.Item(9,5).Value = $v9c5

$v9c6 = @'
//get the first tuple of the first relation using stream.findFirst
  //extract the attribute names using BeanUtils.describe
  //get the first tuple of the second relation using stream.findFirst
  //extract the attribute names using BeanUtils.describe
  //compare the attribute names of both relations and rename if matching exist
  
  //create list<Map<String, Object>>
  //iterate over first relation
   //iterate over second relation
    //create a new map
    //extract values from first relation and put everything into the map
    //do the same with the second relation
    //add map to the list
  
  //return list 
  return l;
  return professorsTemp.stream().collect(Collectors.toList());
 }
'@
$ws.Cells.Item(9,6).Value = $v9c6

$v2c6 = @'
professors.stream().filter(p -> {
    return Integer.valueOf(BeanUtils.getProperty(p, "age").toString()) > 30
      && BeanUtils.getProperty(p, "gender").toString().equals("M");
  }).collect(Collectors.toSet());
'@
$ws.Cells.Item(2,6).Value = $v2c6

$v3c6 = @'
return professorsTemp.stream().map(p -> {
   Map<String, Object> tmp = new HashMap<>();
    tmp.put("name", BeanUtils.getProperty(p, "name"));
    tmp.put("lastName", BeanUtils.getProperty(p, "lastName"));
   return tmp;
  }).collect(Collectors.toSet());
'@
$ws.Cells.Item(3,6).Value = $v3c6

$v4c7 = @'
Can I really simulate renames in Stream?? Why do I need them?
'@
$ws.Cells.Item(4,7).Value = $v4c7

$v4c6 = @'
return professorsTemp.stream().map(p -> {
   String s = (String)BeanUtils.getProperty(p, "name");
    BeanUtils.setProperty(p, "newName", s);
    return p;   
  }).collect(Collectors.toList());
'@
$ws.Cells.Item(4,6).Value = $v4c6

$v6c7 = @'
Forcing using maps for this operation?
For now, we will use map as internal temporary tables.
Our process will work only with maps.
'@
$ws.Cells.Item(6,7).Value = $v6c7

# === Row height changes (auto-computed heights after content changes) ===
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 285

# === Rich text formatting on G6 ("Forcing using maps..." + red warning) ===
# Register the red font in the workbook's style table by applying then clearing it on a scratch cell,
# matching how Excel keeps a font-table entry for a color used in an in-cell rich-text run.
$scratch = $ws.Cells.Item(30,30)
$scratch.Value = "x"
$scratch.Font.Color = 255
$scratch.Clear()

$g6 = $ws.Cells.Item(6,7)
$redChars = $g6.Characters(40, 92)
$redChars.Font.Color = 255

# === Update active selection to match final state ===
$ws.Range("G9").Select()
